## Applies:
##  1) Slide 6 table's tableStyleId change
##       {D9DD7296-C82A-4D70-8BFC-29B6F9499A42} -> {47527566-9C5A-49E7-98FA-C336C0101DA6}
##  2) Re-colour the deck's theme ("Integral" colours -> stock "Office" colours)
##       i.e. what a user sees when switching the Design/Theme colour scheme
##       back to the default Office palette (dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{47527566-9C5A-49E7-98FA-C336C0101DA6}")
    }
}

# --- 2) Theme colours -------------------------------------------------------
# ThemeColorScheme slots are ordered: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink (1-based). RGB is set using the usual Win32 BGR packing,
# i.e. 0xBBGGRR, same as VBA's RGB(r,g,b) function.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
